# api excel data input
# Updates a handful of contact e-mail addresses, refreshes the per-sheet
# AutoFilter "_FilterDatabase" defined-name bookkeeping, repositions the
# remembered cell selection on a few sheets, and nudges a couple of column
# widths - mirroring a normal interactive data-entry session in Excel.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Cont adminstrator" sheet (localSheetId 0)
# ---------------------------------------------------------------------
$wsCont = $wb.Worksheets.Item("Cont adminstrator")
$wsCont.Range("C15").Value = "littlepiglesswt811@automation.33mail.com"

# ---------------------------------------------------------------------
# "Receptie" sheet
# ---------------------------------------------------------------------
$wsReceptie = $wb.Worksheets.Item("Receptie")
$wsReceptie.Range("B2").Value = "lilanna3311@staffcalendis.33mail.com"
$wsReceptie.Range("B3").Value = "moraritza2711@staffcalendis.33mail.com"
$wsReceptie.Range("B4").Value = "rimmelplus3311@staffcalendis.33mail.com"
$wsReceptie.Columns.Item(1).ColumnWidth = 38.833333
$wsReceptie.Columns.Item(2).ColumnWidth = 72.333333

# ---------------------------------------------------------------------
# "Locatii" sheet - column width nudge only
# ---------------------------------------------------------------------
$wsLocatii = $wb.Worksheets.Item("Locatii")
$wsLocatii.Columns.Item(5).ColumnWidth = 28.666667

# ---------------------------------------------------------------------
# "Domenii" sheet (localSheetId 3) - column width nudge
# ---------------------------------------------------------------------
$wsDomenii = $wb.Worksheets.Item("Domenii")
$wsDomenii.Columns.Item(2).ColumnWidth = 33.833333

# ---------------------------------------------------------------------
# "Domenii existente" sheet (localSheetId 4) - no direct edits besides
# the refreshed AutoFilter bookkeeping below
# ---------------------------------------------------------------------
$wsDomeniiExistente = $wb.Worksheets.Item("Domenii existente")

# ---------------------------------------------------------------------
# "Angajati" sheet
# ---------------------------------------------------------------------
$wsAngajati = $wb.Worksheets.Item("Angajati")
$wsAngajati.Range("B2").Value = "marilenajohhjss1521@staffcalendis.33mail.com"
$wsAngajati.Range("B3").Value = "ideaforkih3971@staffcalendis.33mail.com"
$wsAngajati.Range("B4").Value = "boomsie4s2861@staffcalendis.33mail.com"
$wsAngajati.Range("B5").Value = "ocarinass3051@staffcalendis.33mail.com"
$wsAngajati.Columns.Item(1).ColumnWidth = 40.166667
$wsAngajati.Columns.Item(2).ColumnWidth = 46.666667
$wsAngajati.Columns.Item(3).ColumnWidth = 29.166667

# ---------------------------------------------------------------------
# Re-save the AutoFilter range on the three filtered sheets. Excel keeps
# appending a fresh "_xlnm._FilterDatabase" + "_0..." hidden defined name
# every time the workbook round-trips through a save with an active
# AutoFilter - reproduce that bookkeeping explicitly since each sheet
# already carries a long history of these duplicates.
$longSuffix = ""
for ($i = 0; $i -lt 155; $i++) {
    $longSuffix = $longSuffix + "_0"
}

$nameCont = $wsCont.Names.Add("__tmpFilterDatabaseCont", "='Cont adminstrator'!`$A`$1:`$A`$19")
$nameCont.Name = "_xlnm._FilterDatabase" + $longSuffix

$nameDomenii = $wsDomenii.Names.Add("__tmpFilterDatabaseDomenii", "=Domenii!`$A`$4:`$A`$7")
$nameDomenii.Name = "_xlnm._FilterDatabase" + $longSuffix

$nameDomeniiExistente = $wsDomeniiExistente.Names.Add("__tmpFilterDatabaseDomeniiExistente", "='Domenii existente'!`$A`$1:`$Q`$15")
$nameDomeniiExistente.Name = "_xlnm._FilterDatabase" + $longSuffix

# ---------------------------------------------------------------------
# Selections - set these last (in tab order) so the final workbook-level
# active sheet/tab ends up back on "Receptie", matching the saved file.
# ---------------------------------------------------------------------
$wsCont.Range("C15").Select()
$wsAngajati.Range("B6").Select()
$wsReceptie.Range("B40").Select()
